$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates for 2020-09-02 data (nombre_aides in column C, montant_total in column D)
$ws.Cells.Item(2, 3).Value = "'898"
$ws.Cells.Item(2, 4).Value = "'2003687.79"
$ws.Cells.Item(4, 3).Value = "'1037"
$ws.Cells.Item(4, 4).Value = "'3691879.47"
$ws.Cells.Item(6, 3).Value = "'677"
$ws.Cells.Item(6, 4).Value = "'2208707.78"
$ws.Cells.Item(14, 3).Value = "'225"
$ws.Cells.Item(14, 4).Value = "'606362.00"
$ws.Cells.Item(16, 3).Value = "'504"
$ws.Cells.Item(16, 4).Value = "'1873774.13"
$ws.Cells.Item(19, 3).Value = "'9"
$ws.Cells.Item(19, 4).Value = "'30418.77"
$ws.Cells.Item(21, 3).Value = "'339"
$ws.Cells.Item(21, 4).Value = "'1205741.00"
$ws.Cells.Item(25, 3).Value = "'121"
$ws.Cells.Item(25, 4).Value = "'291025.13"
$ws.Cells.Item(26, 3).Value = "'138"
$ws.Cells.Item(26, 4).Value = "'458516.96"
$ws.Cells.Item(27, 3).Value = "'115"
$ws.Cells.Item(27, 4).Value = "'365094.60"
$ws.Cells.Item(31, 3).Value = "'587"
$ws.Cells.Item(31, 4).Value = "'2476919.89"
$ws.Cells.Item(33, 3).Value = "'404"
$ws.Cells.Item(33, 4).Value = "'1422480.57"
$ws.Cells.Item(36, 3).Value = "'372"
$ws.Cells.Item(36, 4).Value = "'891760.71"
$ws.Cells.Item(37, 3).Value = "'221"
$ws.Cells.Item(37, 4).Value = "'693504.04"
$ws.Cells.Item(38, 3).Value = "'210"
$ws.Cells.Item(38, 4).Value = "'575671.14"
$ws.Cells.Item(39, 3).Value = "'7"
$ws.Cells.Item(39, 4).Value = "'17000.00"
$ws.Cells.Item(46, 3).Value = "'415"
$ws.Cells.Item(46, 4).Value = "'1139154.43"
$ws.Cells.Item(48, 3).Value = "'649"
$ws.Cells.Item(48, 4).Value = "'2668879.99"
$ws.Cells.Item(49, 3).Value = "'447"
$ws.Cells.Item(49, 4).Value = "'1605996.50"
$ws.Cells.Item(52, 3).Value = "'3824"
$ws.Cells.Item(52, 4).Value = "'8816141.47"
$ws.Cells.Item(55, 3).Value = "'4057"
$ws.Cells.Item(55, 4).Value = "'14289025.38"
$ws.Cells.Item(60, 3).Value = "'83"
$ws.Cells.Item(60, 4).Value = "'193416.69"
$ws.Cells.Item(61, 3).Value = "'141"
$ws.Cells.Item(61, 4).Value = "'415368.00"
$ws.Cells.Item(62, 3).Value = "'41"
$ws.Cells.Item(62, 4).Value = "'135177.00"
$ws.Cells.Item(75, 3).Value = "'395"
$ws.Cells.Item(75, 4).Value = "'1004406.70"
$ws.Cells.Item(77, 3).Value = "'945"
$ws.Cells.Item(77, 4).Value = "'3328064.26"
$ws.Cells.Item(78, 3).Value = "'530"
$ws.Cells.Item(78, 4).Value = "'1751225.47"
